# Figure 3: panel tags to uppercase
#
# Slide 1 has a group ("Group 33") containing two small label rectangles
# that tag the two panels of the figure: shape id=44 ("Rectangle 43")
# holds the text "(a)" and shape id=45 ("Rectangle 44") holds "(b)".
# Uppercase both panel tags: "(a)" -> "(A)" and "(b)" -> "(B)".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-Shape {
    param($shapes, [int]$id, [string]$name, [string]$text)
    $byName = $null
    $byText = $null
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)

        if ($shp.Id -eq $id) {
            return $shp
        }
        if ($name -and $shp.Name -eq $name) {
            $byName = $shp
        }
        if ($text -and $shp.HasTextFrame -and ($shp.TextFrame.TextRange.Text -eq $text)) {
            $byText = $shp
        }

        if ($shp.Type -eq 6) {
            $found = Find-Shape $shp.GroupItems $id $name $text
            if ($found) { return $found }
        }
    }
    if ($byName) { return $byName }
    if ($byText) { return $byText }
    return $null
}

$rectA = Find-Shape $s.Shapes 44 "Rectangle 43" "(a)"
$rectB = Find-Shape $s.Shapes 45 "Rectangle 44" "(b)"

$rectA.TextFrame.TextRange.Text = "(A)"
$rectB.TextFrame.TextRange.Text = "(B)"
